$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.3032696666666667
$ws.Cells.Item(2, 8).Value = 0.909809
$ws.Cells.Item(2, 9).Value = 0.150143210583919
$ws.Cells.Item(2, 10).Value = 0.150143210583919
$ws.Cells.Item(2, 13).Value = 82.43338033333333
$ws.Cells.Item(2, 14).Value = 247.300141
$ws.Cells.Item(2, 15).Value = 0.3670006993429558
$ws.Cells.Item(2, 16).Value = 0.3670006993429557
$ws.Cells.Item(2, 17).Value = 24.99954377589656
$ws.Cells.Item(2, 18).Value = 224.995893983069
$ws.Cells.Item(2, 19).Value = 0.05510266328589494
$ws.Cells.Item(2, 20).Value = 0.05510266328589493

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.3032696666666667
$ws.Cells.Item(3, 8).Value = 0.909809
$ws.Cells.Item(3, 9).Value = 0.150143210583919
$ws.Cells.Item(3, 10).Value = 0.150143210583919
$ws.Cells.Item(3, 15).Value = 0.3956886215996139
$ws.Cells.Item(3, 16).Value = 0.3956886215996139
$ws.Cells.Item(3, 17).Value = 26.95372252699655
$ws.Cells.Item(3, 18).Value = 242.583502742969
$ws.Cells.Item(3, 19).Value = 0.05940996003849146
$ws.Cells.Item(3, 20).Value = 0.05940996003849145

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.3032696666666667
$ws.Cells.Item(4, 8).Value = 0.909809
$ws.Cells.Item(4, 9).Value = 0.150143210583919
$ws.Cells.Item(4, 10).Value = 0.150143210583919
$ws.Cells.Item(4, 13).Value = 42.93483766666667
$ws.Cells.Item(4, 14).Value = 128.804513
$ws.Cells.Item(4, 15).Value = 0.1911496942879982
$ws.Cells.Item(4, 16).Value = 0.1911496942879981
$ws.Cells.Item(4, 17).Value = 13.02083390755745
$ws.Cells.Item(4, 18).Value = 117.187505168017
$ws.Cells.Item(4, 19).Value = 0.02869982880253464
$ws.Cells.Item(4, 20).Value = 0.02869982880253464

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.3032696666666667
$ws.Cells.Item(5, 8).Value = 0.909809
$ws.Cells.Item(5, 9).Value = 0.150143210583919
$ws.Cells.Item(5, 10).Value = 0.150143210583919
$ws.Cells.Item(5, 13).Value = 10.368389
$ws.Cells.Item(5, 14).Value = 31.105167
$ws.Cells.Item(5, 15).Value = 0.04616098476943217
$ws.Cells.Item(5, 16).Value = 0.04616098476943217
$ws.Cells.Item(5, 17).Value = 3.144417875900333
$ws.Cells.Item(5, 18).Value = 28.299760883103
$ws.Cells.Item(5, 19).Value = 0.006930758456997931
$ws.Cells.Item(5, 20).Value = 0.006930758456997929

$ws.Cells.Item(6, 9).Value = 0.6441382016790526
$ws.Cells.Item(6, 10).Value = 0.6441382016790526
$ws.Cells.Item(6, 13).Value = 82.43338033333333
$ws.Cells.Item(6, 14).Value = 247.300141
$ws.Cells.Item(6, 15).Value = 0.3670006993429558
$ws.Cells.Item(6, 16).Value = 0.3670006993429557
$ws.Cells.Item(6, 17).Value = 107.2520103171917
$ws.Cells.Item(6, 18).Value = 965.268092854725
$ws.Cells.Item(6, 19).Value = 0.2363991704897262
$ws.Cells.Item(6, 20).Value = 0.2363991704897261

$ws.Cells.Item(7, 9).Value = 0.6441382016790526
$ws.Cells.Item(7, 10).Value = 0.6441382016790526
$ws.Cells.Item(7, 15).Value = 0.3956886215996139
$ws.Cells.Item(7, 16).Value = 0.3956886215996139
$ws.Cells.Item(7, 19).Value = 0.2548781571420384
$ws.Cells.Item(7, 20).Value = 0.2548781571420384

$ws.Cells.Item(8, 9).Value = 0.6441382016790526
$ws.Cells.Item(8, 10).Value = 0.6441382016790526
$ws.Cells.Item(8, 13).Value = 42.93483766666667
$ws.Cells.Item(8, 14).Value = 128.804513
$ws.Cells.Item(8, 15).Value = 0.1911496942879982
$ws.Cells.Item(8, 16).Value = 0.1911496942879981
$ws.Cells.Item(8, 17).Value = 55.86144391715834
$ws.Cells.Item(8, 18).Value = 502.7529952544251
$ws.Cells.Item(8, 19).Value = 0.1231268203301718
$ws.Cells.Item(8, 20).Value = 0.1231268203301718

$ws.Cells.Item(9, 9).Value = 0.6441382016790526
$ws.Cells.Item(9, 10).Value = 0.6441382016790526
$ws.Cells.Item(9, 13).Value = 10.368389
$ws.Cells.Item(9, 14).Value = 31.105167
$ws.Cells.Item(9, 15).Value = 0.04616098476943217
$ws.Cells.Item(9, 16).Value = 0.04616098476943217
$ws.Cells.Item(9, 17).Value = 13.490051718175
$ws.Cells.Item(9, 18).Value = 121.410465463575
$ws.Cells.Item(9, 19).Value = 0.02973405371711618
$ws.Cells.Item(9, 20).Value = 0.02973405371711617

$ws.Cells.Item(10, 7).Value = 0.4155246666666666
$ws.Cells.Item(10, 8).Value = 1.246574
$ws.Cells.Item(10, 9).Value = 0.2057185877370285
$ws.Cells.Item(10, 10).Value = 0.2057185877370285
$ws.Cells.Item(10, 13).Value = 82.43338033333333
$ws.Cells.Item(10, 14).Value = 247.300141
$ws.Cells.Item(10, 15).Value = 0.3670006993429558
$ws.Cells.Item(10, 16).Value = 0.3670006993429557
$ws.Cells.Item(10, 17).Value = 34.25310288521488
$ws.Cells.Item(10, 18).Value = 308.2779259669339
$ws.Cells.Item(10, 19).Value = 0.07549886556733466
$ws.Cells.Item(10, 20).Value = 0.07549886556733465

$ws.Cells.Item(11, 7).Value = 0.4155246666666666
$ws.Cells.Item(11, 8).Value = 1.246574
$ws.Cells.Item(11, 9).Value = 0.2057185877370285
$ws.Cells.Item(11, 10).Value = 0.2057185877370285
$ws.Cells.Item(11, 15).Value = 0.3956886215996139
$ws.Cells.Item(11, 16).Value = 0.3956886215996139
$ws.Cells.Item(11, 17).Value = 36.93061917981488
$ws.Cells.Item(11, 18).Value = 332.3755726183339
$ws.Cells.Item(11, 19).Value = 0.08140050441908404
$ws.Cells.Item(11, 20).Value = 0.08140050441908404

$ws.Cells.Item(12, 7).Value = 0.4155246666666666
$ws.Cells.Item(12, 8).Value = 1.246574
$ws.Cells.Item(12, 9).Value = 0.2057185877370285
$ws.Cells.Item(12, 10).Value = 0.2057185877370285
$ws.Cells.Item(12, 13).Value = 42.93483766666667
$ws.Cells.Item(12, 14).Value = 128.804513
$ws.Cells.Item(12, 15).Value = 0.1911496942879982
$ws.Cells.Item(12, 16).Value = 0.1911496942879981
$ws.Cells.Item(12, 17).Value = 17.84048410982911
$ws.Cells.Item(12, 18).Value = 160.564356988462
$ws.Cells.Item(12, 19).Value = 0.03932304515529172
$ws.Cells.Item(12, 20).Value = 0.03932304515529172

$ws.Cells.Item(13, 7).Value = 0.4155246666666666
$ws.Cells.Item(13, 8).Value = 1.246574
$ws.Cells.Item(13, 9).Value = 0.2057185877370285
$ws.Cells.Item(13, 10).Value = 0.2057185877370285
$ws.Cells.Item(13, 13).Value = 10.368389
$ws.Cells.Item(13, 14).Value = 31.105167
$ws.Cells.Item(13, 15).Value = 0.04616098476943217
$ws.Cells.Item(13, 16).Value = 0.04616098476943217
$ws.Cells.Item(13, 17).Value = 4.308321383095333
$ws.Cells.Item(13, 18).Value = 38.77489244785799
$ws.Cells.Item(13, 19).Value = 0.009496172595318068
$ws.Cells.Item(13, 20).Value = 0.009496172595318068
